$wb = $excel.ActiveWorkbook

$wsSpreadsheet = $wb.Worksheets.Item("TABLE_FROM_SPREADSHEET")
$wsSpreadsheet.Range("B1").Value = "COLUMN_A"

$wsCsv = $wb.Worksheets.Item("TABLE_FROM_CSV")
$wsCsv.Range("B1").Value = "COLUMN_B"

$wsDatetime = $wb.Worksheets.Item("TABLE_DATETIME")
$wsDatetime.Range("G1").Value = "TOUCHED UP"
$wsDatetime.Range("H1").Value = "TOUCHED UP_raw"
